# Add a new "Player Info" worksheet as the first sheet in the workbook,
# rename the MATCH_CARD_LINK columns to MATCH_CODE on the existing
# "ODI Batting"/"ODI Bowling" sheets, and replace the full scorecard URLs
# in those columns with just the bare numeric match code.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Player Info" sheet before the first sheet -------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row
$playerInfo.Range("A2").Value = "'5971"
$playerInfo.Range("B2").Value = "Venkatesh Rajasekaran Iyer"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium"

$playerInfo.Range("A1").Select()

# Re-fetch the other sheets by name now that sheet order has changed.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# --- 2. ODI Batting: rename MATCH_CARD_LINK -> MATCH_CODE, trim URLs ----
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").Value = "'4524"
$battingSheet.Range("D3").Value = "'4526"

# --- 3. ODI Bowling: rename MATCH_CARD_LINK -> MATCH_CODE, trim URLs ----
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").Value = "'4526"
